$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text values that are unambiguous (not parsed as pure numbers)
$ws.Range("D2").Value = "26.423.89"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "1.839.44"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  -6.64%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  -1.96%  "
$ws.Range("E8").Value = "  -5.20%  "
$ws.Range("E9").Value = "  -2.68%  "
$ws.Range("E10").Value = "  -7.60%  "
$ws.Range("E11").Value = "  -4.93%  "
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").Value = "1.837.45"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("E15").Value = "  -2.94%  "
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("E17").Value = "  -4.72%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").Value = "26.442.15"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").Value = "2.070.68"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("E22").Value = "  -4.09%  "
$ws.Range("E23").Value = "  -5.73%  "
$ws.Range("E24").Value = "  -4.24%  "
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("E26").Value = "  -6.38%  "
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("E30").Value = "  -4.38%  "
$ws.Range("E31").Value = "  -4.52%  "
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("E36").Value = "  -3.52%  "
$ws.Range("E37").Value = "  -7.04%  "
$ws.Range("E38").Value = "  -5.39%  "
$ws.Range("E39").Value = "  -7.05%  "
$ws.Range("E40").Value = "  -5.87%  "
$ws.Range("E41").Value = "  -4.39%  "
$ws.Range("E42").Value = "  -7.27%  "
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("E45").Value = "  -5.53%  "
$ws.Range("E47").Value = "  -8.02%  "
$ws.Range("E48").Value = "  -3.23%  "
$ws.Range("E50").Value = "  -9.32%  "
$ws.Range("E51").Value = "  -0.07%  "

# Set numeric-looking text values, forcing text type then restoring default style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5206"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3263"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06775"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7648"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07692"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.023"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007943"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.562"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.461"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.946"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.215"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.652"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.161"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.120"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08706"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04786"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.124"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.841"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7013"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.059"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01750"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.191"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4818"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "110.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8883"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.658"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05869"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4110"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.996"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1217"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8847"
$ws.Range("D51").Style = "Normal"
